# Rename the four auto-generated TOC bookmarks that Word regenerates whenever
# the document's table of contents is rebuilt. Each `_Toc4792545xx` bookmark
# becomes `_Toc4834968xx`, keeping it anchored to exactly the same text span.
#
# The Word OM has no direct "rename" for a Bookmark, so each one is recreated:
# remember its Range, delete it, then add a new bookmark with the new name at
# that same Range.

$d = $word.ActiveDocument

$renames = @(
    @{ Old = "_Toc479254566"; New = "_Toc483496845" },
    @{ Old = "_Toc479254567"; New = "_Toc483496846" },
    @{ Old = "_Toc479254568"; New = "_Toc483496847" },
    @{ Old = "_Toc479254569"; New = "_Toc483496848" }
)

foreach ($pair in $renames) {
    $bm = $d.Bookmarks($pair.Old)
    $r = $bm.Range
    $bm.Delete()
    $d.Bookmarks.Add($pair.New, $r)
}
